# UC001 - Alocar Professor em uma Disciplina
# Commit: "2# fechamento da primeira entraga / o termo dupla so se for de 1"
#
# Changes applied:
#  - C4  (Fluxo Principal, Resultado esperado #1): reworded text about the
#    new allocation screen (now mentions two separate lists).
#  - B11 (Fluxo Alternativo, Passos #1): reworded the trigger condition from
#    "disponíveis" to "cadastrados".
#  - C11 (Fluxo Alternativo, Resultado esperado #1): reworded the warning
#    message shown to the administrator.
#  - Row 4 height grew (explicit resize) and row 11 height grew (more text).
#  - Selected cell moved from C13 to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c4Text = 'Uma nova tela é exibida contendo uma lista dos professores e outra lista com as disciplinas que ainda não tem professores alocados.'
$b11Text = 'Clicar na opção "Alocar Professor em uma disciplina" sem haver professor(es) e/ou disciplina(s) cadastrados.'
$c11Text = 'É Exibida na tela uma mensagem, informando  para o administrador que não há  "Não a Disciplinas Cadastrada!! "para caso da disciplina disciplina, e para cso do professor "Não ha professor Cadastrado" ' + "`n`n" + 'A operação é encerrada.'

$ws.Range("C4").Value = $c4Text
$ws.Range("B11").Value = $b11Text
$ws.Range("C11").Value = $c11Text

# Row height adjustments
$ws.Rows(4).RowHeight = 28.5
$ws.Rows(11).RowHeight = 63.75

# Move the active selection to B13 (was C13)
$null = $ws.Range("B13").Select()
